# Add district heat / hydrogen carbon tax
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("BCTR")

# Row 7 = "district heat and hydrogen sector": B7 = B6*0.75*0.6, C7:AE7 copy the
# same relative formula across the row (matching si="3" shared formula in the diff).
$ws.Range("B7:AE7").Formula = "=B6*0.75*0.6"

# Update the selection shown in the sheet view to match the new active cell/range.
$ws.Range("B7:AE7").Select()

$wb.Save()
